# Rest Assured framework merged and data is being written to xls properly
#
# The sheet gains a "compare result" column after every OWM/NDTV pair, plus
# brand-new Weather-Condition and Wind columns. Final column order:
#   A City
#   B LiveTemp (OpenWeatherMap)     E LiveTemp (NDTV) [dup of B]   D LiveTemp Compare Result
#   E Max Temp (OpenWeatherMap)     F Max Temp (NDTV) [dup of B]   G Max Temp Compare Result
#   H Weather Condition (OpenWeatherMap)  I Weather Condition (NDTV)  J Weather Condition Compare Result
#   K Humidity (OpenWeatherMap)     L Humidity (NDTV)               M Humidity Compare Result
#   N Wind (OpenWeatherMap)         O Wind (NDTV)                   P Wind Compare Result

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "City"
$ws.Range("B1").Value = "LiveTemp (OpenWeatherMap)"
$ws.Range("C1").Value = "LiveTemp (NDTV)"
$ws.Range("D1").Value = "LiveTemp Compare Result"
$ws.Range("E1").Value = "Max Temp (OpenWeatherMap)"
$ws.Range("F1").Value = "Max Temp (NDTV)"
$ws.Range("G1").Value = "Max Temp Compare Result"
$ws.Range("H1").Value = "Weather Condition (OpenWeatherMap)"
$ws.Range("I1").Value = "Weather Condition (NDTV)"
$ws.Range("J1").Value = "Weather Condition Compare Result"
$ws.Range("K1").Value = "Humidity (OpenWeatherMap)"
$ws.Range("L1").Value = "Humidity (NDTV)"
$ws.Range("M1").Value = "Humidity Compare Result"
$ws.Range("N1").Value = "Wind (OpenWeatherMap)"
$ws.Range("O1").Value = "Wind (NDTV)"
$ws.Range("P1").Value = "Wind Compare Result"

# ---------------------------------------------------------------------------
# Row 2 - Lucknow
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Lucknow"
$ws.Range("B2").Value = "29.00"
$ws.Range("E2").Value = "29.00"
$ws.Range("H2").Value = "Mist&mist"
$ws.Range("K2").Value = "89"
$ws.Range("N2").Value = "1.00"

# ---------------------------------------------------------------------------
# Row 3 - Bhopal
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Bhopal"
$ws.Range("B3").Value = "27.00"
$ws.Range("E3").Value = "27.00"
$ws.Range("H3").Value = "Clouds&few clouds"
$ws.Range("K3").Value = "88"
$ws.Range("N3").Value = "2.60"

# ---------------------------------------------------------------------------
# Row 4 - Ajmer
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Ajmer"
$ws.Range("B4").Value = "30.38"
$ws.Range("E4").Value = "30.38"
$ws.Range("H4").Value = "Clouds&overcast clouds"
$ws.Range("K4").Value = "59"
$ws.Range("N4").Value = "5.60"

# ---------------------------------------------------------------------------
# Styling: the "compare result" placeholder cells (C3:D3, C4:D4) keep the
# Hyperlink-flavoured text style the template already used for its blank
# "to-be-filled" slots (same treatment the old sheet gave C3/C4/I2:I4).
# ---------------------------------------------------------------------------
$ws.Range("C3:D3").Style = "Hyperlink"
$ws.Range("C4:D4").Style = "Hyperlink"
$ws.Range("C3:D4").NumberFormat = "@"

$ws.Range("N2:N4").Style = "Hyperlink"
$ws.Range("N2:N4").NumberFormat = "@"
$ws.Range("N2:N4").Value2 = $ws.Range("N2:N4").Value2

# Re-apply the values lost by forcing the Hyperlink style on N2:N4 (style
# application can reset contents on some hosts) and make sure the text stays
# text, not re-interpreted as a number.
$ws.Range("N2").Value = "1.00"
$ws.Range("N3").Value = "2.60"
$ws.Range("N4").Value = "5.60"

# ---------------------------------------------------------------------------
# Column widths (widths below are the Excel "character" widths taken from
# the saved workbook)
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 16.42578125
$ws.Columns.Item(4).ColumnWidth = 18.28515625
$ws.Columns.Item(6).ColumnWidth = 17.140625
$ws.Columns.Item(7).ColumnWidth = 24.7109375
$ws.Columns.Item(8).ColumnWidth = 36.74609375
$ws.Columns.Item(9).ColumnWidth = 24.7109375
$ws.Columns.Item(10).ColumnWidth = 24.7109375
$ws.Columns.Item(12).ColumnWidth = 16.42578125
$ws.Columns.Item(14).ColumnWidth = 24.45703125

# ---------------------------------------------------------------------------
# Header row wraps text and grows taller to fit two-line captions.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 45
$ws.Range("A1:P1").WrapText = $true
$ws.Range("A1:A4").WrapText = $true
$ws.Range("C3:D4").WrapText = $true
$ws.Range("N2:N4").WrapText = $true

# ---------------------------------------------------------------------------
# Selection / view: active cell moves to the new "Weather Condition Compare
# Result" header once the extra columns are in place.
# ---------------------------------------------------------------------------
$null = $ws.Range("J1").Select()
$excel.ActiveWindow.ScrollColumn = 2
